$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3.0
$ws.Range("F2").Value = 1.0
$ws.Range("G2").Value = 0.350653
$ws.Range("H2").Value = 1.051959
$ws.Range("I2").Value = 0.04536179359243143
$ws.Range("J2").Value = 0.04536179359243143
$ws.Range("M2").Value = 0.789222
$ws.Range("N2").Value = 2.367666
$ws.Range("O2").Value = 0.01341929863527565
$ws.Range("P2").Value = 0.01341929863527565
$ws.Range("Q2").Value = 0.276743061966
$ws.Range("R2").Value = 2.490687557694
$ws.Range("S2").Value = 0.0006087234548485707
$ws.Range("T2").Value = 0.0006087234548485707
$ws.Range("E3").Value = 3.0
$ws.Range("F3").Value = 1.0
$ws.Range("G3").Value = 0.350653
$ws.Range("H3").Value = 1.051959
$ws.Range("I3").Value = 0.04536179359243143
$ws.Range("J3").Value = 0.04536179359243143
$ws.Range("O3").Value = 0.005047365584441773
$ws.Range("P3").Value = 0.005047365584441773
$ws.Range("Q3").Value = 0.104090641744
$ws.Range("R3").Value = 0.936815775696
$ws.Range("S3").Value = 0.0002289575558269898
$ws.Range("T3").Value = 0.0002289575558269898
$ws.Range("E4").Value = 3.0
$ws.Range("F4").Value = 1.0
$ws.Range("G4").Value = 0.350653
$ws.Range("H4").Value = 1.051959
$ws.Range("I4").Value = 0.04536179359243143
$ws.Range("J4").Value = 0.04536179359243143
$ws.Range("M4").Value = 57.61405833333333
$ws.Range("N4").Value = 172.842175
$ws.Range("O4").Value = 0.9796232927683105
$ws.Range("P4").Value = 0.9796232927683105
$ws.Range("Q4").Value = 20.20254239675834
$ws.Range("R4").Value = 181.822881570825
$ws.Range("S4").Value = 0.04443746960489413
$ws.Range("T4").Value = 0.04443746960489413
$ws.Range("E5").Value = 3.0
$ws.Range("F5").Value = 1.0
$ws.Range("G5").Value = 0.350653
$ws.Range("H5").Value = 1.051959
$ws.Range("I5").Value = 0.04536179359243143
$ws.Range("J5").Value = 0.04536179359243143
$ws.Range("M5").Value = 0.1123343333333333
$ws.Range("N5").Value = 0.337003
$ws.Range("O5").Value = 0.001910043011972043
$ws.Range("P5").Value = 0.001910043011972043
$ws.Range("Q5").Value = 0.03939037098633334
$ws.Range("R5").Value = 0.354513338877
$ws.Range("S5").Value = 0.00008664297686174185
$ws.Range("T5").Value = 0.00008664297686174186
$ws.Range("I6").Value = 0.8482855786262421
$ws.Range("J6").Value = 0.8482855786262421
$ws.Range("M6").Value = 0.789222
$ws.Range("N6").Value = 2.367666
$ws.Range("O6").Value = 0.01341929863527565
$ws.Range("P6").Value = 0.01341929863527565
$ws.Range("Q6").Value = 5.175217509251999
$ws.Range("R6").Value = 46.57695758326799
$ws.Range("S6").Value = 0.01138339750758314
$ws.Range("T6").Value = 0.01138339750758314
$ws.Range("I7").Value = 0.8482855786262421
$ws.Range("J7").Value = 0.8482855786262421
$ws.Range("O7").Value = 0.005047365584441773
$ws.Range("P7").Value = 0.005047365584441773
$ws.Range("S7").Value = 0.00428160743533637
$ws.Range("T7").Value = 0.00428160743533637
$ws.Range("I8").Value = 0.8482855786262421
$ws.Range("J8").Value = 0.8482855786262421
$ws.Range("M8").Value = 57.61405833333333
$ws.Range("N8").Value = 172.842175
$ws.Range("O8").Value = 0.9796232927683105
$ws.Range("P8").Value = 0.9796232927683105
$ws.Range("Q8").Value = 377.7964672370166
$ws.Range("R8").Value = 3400.16820513315
$ws.Range("S8").Value = 0.8310003117417109
$ws.Range("T8").Value = 0.8310003117417109
$ws.Range("I9").Value = 0.8482855786262421
$ws.Range("J9").Value = 0.8482855786262421
$ws.Range("M9").Value = 0.1123343333333333
$ws.Range("N9").Value = 0.337003
$ws.Range("O9").Value = 0.001910043011972043
$ws.Range("P9").Value = 0.001910043011972043
$ws.Range("Q9").Value = 0.7366173380326665
$ws.Range("R9").Value = 6.629556042293999
$ws.Range("S9").Value = 0.001620261941611715
$ws.Range("T9").Value = 0.001620261941611715
$ws.Range("G10").Value = 0.7457606666666666
$ws.Range("H10").Value = 2.237282
$ws.Range("I10").Value = 0.09647441040198541
$ws.Range("J10").Value = 0.09647441040198541
$ws.Range("M10").Value = 0.789222
$ws.Range("N10").Value = 2.367666
$ws.Range("O10").Value = 0.01341929863527565
$ws.Range("P10").Value = 0.01341929863527565
$ws.Range("Q10").Value = 0.588570724868
$ws.Range("R10").Value = 5.297136523811999
$ws.Range("S10").Value = 0.001294618923846386
$ws.Range("T10").Value = 0.001294618923846386
$ws.Range("G11").Value = 0.7457606666666666
$ws.Range("H11").Value = 2.237282
$ws.Range("I11").Value = 0.09647441040198541
$ws.Range("J11").Value = 0.09647441040198541
$ws.Range("O11").Value = 0.005047365584441773
$ws.Range("P11").Value = 0.005047365584441773
$ws.Range("Q11").Value = 0.2213775623786667
$ws.Range("R11").Value = 1.992398061408
$ws.Range("S11").Value = 0.0004869416188422926
$ws.Range("T11").Value = 0.0004869416188422926
$ws.Range("G12").Value = 0.7457606666666666
$ws.Range("H12").Value = 2.237282
$ws.Range("I12").Value = 0.09647441040198541
$ws.Range("J12").Value = 0.09647441040198541
$ws.Range("M12").Value = 57.61405833333333
$ws.Range("N12").Value = 172.842175
$ws.Range("O12").Value = 0.9796232927683105
$ws.Range("P12").Value = 0.9796232927683105
$ws.Range("Q12").Value = 42.96629855203889
$ws.Range("R12").Value = 386.69668696835
$ws.Range("S12").Value = 0.0945085795858743
$ws.Range("T12").Value = 0.0945085795858743
$ws.Range("G13").Value = 0.7457606666666666
$ws.Range("H13").Value = 2.237282
$ws.Range("I13").Value = 0.09647441040198541
$ws.Range("J13").Value = 0.09647441040198541
$ws.Range("M13").Value = 0.1123343333333333
$ws.Range("N13").Value = 0.337003
$ws.Range("O13").Value = 0.001910043011972043
$ws.Range("P13").Value = 0.001910043011972043
$ws.Range("Q13").Value = 0.08377452731622222
$ws.Range("R13").Value = 0.753970745846
$ws.Range("S13").Value = 0.0001842702734224352
$ws.Range("T13").Value = 0.0001842702734224352
$ws.Range("G14").Value = 0.07636
$ws.Range("H14").Value = 0.22908
$ws.Range("I14").Value = 0.009878217379341012
$ws.Range("J14").Value = 0.009878217379341012
$ws.Range("M14").Value = 0.789222
$ws.Range("N14").Value = 2.367666
$ws.Range("O14").Value = 0.01341929863527565
$ws.Range("P14").Value = 0.01341929863527565
$ws.Range("Q14").Value = 0.06026499192
$ws.Range("R14").Value = 0.54238492728
$ws.Range("S14").Value = 0.000132558748997547
$ws.Range("T14").Value = 0.000132558748997547
$ws.Range("G15").Value = 0.07636
$ws.Range("H15").Value = 0.22908
$ws.Range("I15").Value = 0.009878217379341012
$ws.Range("J15").Value = 0.009878217379341012
$ws.Range("O15").Value = 0.005047365584441773
$ws.Range("P15").Value = 0.005047365584441773
$ws.Range("Q15").Value = 0.02266731328
$ws.Range("R15").Value = 0.20400581952
$ws.Range("S15").Value = 0.00004985897443612043
$ws.Range("T15").Value = 0.00004985897443612043
$ws.Range("G16").Value = 0.07636
$ws.Range("H16").Value = 0.22908
$ws.Range("I16").Value = 0.009878217379341012
$ws.Range("J16").Value = 0.009878217379341012
$ws.Range("M16").Value = 57.61405833333333
$ws.Range("N16").Value = 172.842175
$ws.Range("O16").Value = 0.9796232927683105
$ws.Range("P16").Value = 0.9796232927683105
$ws.Range("Q16").Value = 4.399409494333333
$ws.Range("R16").Value = 39.594685449
$ws.Range("S16").Value = 0.009676931835831192
$ws.Range("T16").Value = 0.009676931835831192
$ws.Range("G17").Value = 0.07636
$ws.Range("H17").Value = 0.22908
$ws.Range("I17").Value = 0.009878217379341012
$ws.Range("J17").Value = 0.009878217379341012
$ws.Range("M17").Value = 0.1123343333333333
$ws.Range("N17").Value = 0.337003
$ws.Range("O17").Value = 0.001910043011972043
$ws.Range("P17").Value = 0.001910043011972043
$ws.Range("Q17").Value = 0.008577849693333332
$ws.Range("R17").Value = 0.07720064724
$ws.Range("S17").Value = 0.00001886782007615108
$ws.Range("T17").Value = 0.00001886782007615109
